$wb = $excel.ActiveWorkbook

# --- CoverageOfficer sheet: new officer added in column A, previous
#     name shifted into column B -------------------------------------------
$wsCoverage = $wb.Worksheets.Item("CoverageOfficer")
$wsCoverage.Range("A2").Value = "Michael Morabito"
$wsCoverage.Range("B2").Value = "Jim Lavelle"
$wsCoverage.Range("C2").Clear()

# --- Companies sheet: "Varta AG" record replaced by "ADK Holdings",
#     previous value moved into column D -----------------------------------
$wsCompanies = $wb.Worksheets.Item("Companies")
$wsCompanies.Range("A3").Value = "ADK Holdings"
$wsCompanies.Range("D3").Value = "Varta AG"

# --- CoverageOfficer sheet: second new officer added in column A, previous
#     name shifted into column B -------------------------------------------
$wsCoverage.Range("A3").Value = "Yuta Nakamura"
$wsCoverage.Range("B3").Value = "Steve Hughes"

[void]$wsCoverage.Select()
[void]$wsCoverage.Range("F22").Select()

[void]$wsCompanies.Select()
[void]$wsCompanies.Range("C6").Select()
